# LarvalAbundSensitivity.xlsx revision
# - Updates mean larval abundance / positive-station values with recomputed numbers
# - Adds a new "Larval index" (column E) value for most rows
# - Refreshes the "Notes" column text with corrected stats and caveats
# - Bumps row 6's "Number of days" and row 8 config label text stays the same
# - Moves the active selection to F10

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("C2").Value = 2.79
$ws.Range("D2").Value = 11.75
$ws.Range("E2").Value = 1.24

# --- Row 3 ---
$ws.Range("C3").Value = 2.79
$ws.Range("D3").Value = 13
$ws.Range("E3").Value = 1.29

# --- Row 4 ---
$ws.Range("C4").Value = 2.0499999999999998
$ws.Range("D4").Value = 12.15
$ws.Range("E4").Value = 0.94

# --- Row 5 ---
$ws.Range("C5").Value = 3.19
$ws.Range("D5").Value = 13
$ws.Range("E5").Value = 1.48

# --- Row 6 ---
$ws.Range("B6").Value = 58
$ws.Range("C6").Value = 1.94
$ws.Range("D6").Value = 12.94
$ws.Range("E6").Value = 0.81
$ws.Range("F6").Value = "Area stratified mean includes shallow stations; Offshore mean = 2.16 (12.49), shelfbreak mean = 1.15 (14.57)"

# --- Row 7 ---
$ws.Range("C7").Value = 2.5499999999999998
$ws.Range("D7").Value = 12.94
$ws.Range("E7").Value = 1.1200000000000001
$ws.Range("F7").Value = "Area stratified mean includes shallow stations; Offshore mean = 2.68 (12.49), shelfbreak mean = 2.08 (14.57)"

# --- Row 8 ---
$ws.Range("C8").Value = 2.46
$ws.Range("D8").Value = 11.12
$ws.Range("E8").Value = "NA"
$ws.Range("F8").Value = "Area stratified mean includes shallow stations; Offshore mean = 3.04 (12.49), shelfbreak mean = 0.38 (6.15); Can't calculate larval index in this case because there is only 1 station in the shelfbreak region."

# --- Row 9 ---
$ws.Range("F9").Value = "In the SEAMAP data, 35% of the stations are in water shallower than 1000 m. There are also high catches at some of these shallower stations. I can't calculate larval index because I don't have the lengths of individuals."

# --- Selection moves from G11 to F10 ---
$ws.Range("F10").Select()
